# Append a new data row (row 66) to each of the 4 worksheets, mirroring the
# daily-log pattern already present in the sheet (one new row per day).
# Row 65's date/number style is copied onto the new row's A cell so the
# new row renders with the same YYYY-MM-DD HH:MM:SS format.

$wb = $excel.ActiveWorkbook

# Per-sheet values for the new row 66, in column order: A,B,C,D,E,F,G,H,I
$rows = @{
    "MID_LFT_#1" = @{
        A = 45852.46226851852
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
        D = "0x01,0x54"
        E = "0x07"
        F = 400
        G = [double]"5.68631262647113e+23"
        H = 340
        I = 7
    }
    "MID_LFT_#2" = @{
        A = 45852.46226851852
        B = "0x01,0x7c"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D = "0x01,0x50"
        E = "0x19"
        F = 380
        G = [double]"5.68432987514711e+23"
        H = 336
        I = 25
    }
    "MID_PLT_#1" = @{
        A = 45852.46226851852
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D = "0x00,0x66"
        E = "0x15"
        F = 110
        G = [double]"5.68631262647113e+23"
        H = 102
        I = 15
    }
    "MID_PLT_#2" = @{
        A = 45852.46226851852
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D = "0x00,0x7B"
        E = "0x9"
        F = 130
        G = [double]"5.68631262647113e+23"
        H = 123
        I = 9
    }
}

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $data = $rows[$ws.Name]
    if ($data -eq $null) { continue }

    $newRow = 66

    # Preserve the date/time number format used by the existing rows in
    # column A (style index 2 -> "YYYY-MM-DD HH:MM:SS").
    $dateFormat = $ws.Range("A65").NumberFormat

    $ws.Cells.Item($newRow, 1).Value = $data.A
    $ws.Cells.Item($newRow, 1).NumberFormat = $dateFormat

    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E
    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = $data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}
